$d = $word.ActiveDocument

# --- 1) Date line near top: "Date: 27th January 2024" -> "Date: 1st March 2024"
$p2 = $d.Paragraphs(2).Range
$p2.Find.Execute(" 27", $true, $false, $false, $false, $false, $true, 1, $false, " 1", 2)

$p2 = $d.Paragraphs(2).Range
$p2.Find.Execute("th", $true, $false, $false, $false, $false, $true, 1, $false, "st", 2)

$p2 = $d.Paragraphs(2).Range
$p2.Find.Execute("January", $true, $false, $false, $false, $false, $true, 1, $false, "March", 2)

# --- 2) Name line: "Name: Priyesh Gaude" -> "Name: Ramdas Tuyenkar"
$p3 = $d.Paragraphs(3).Range
$p3.Find.Execute("Priyesh", $true, $false, $false, $false, $false, $true, 1, $false, "Ramdas", 2)

$p3 = $d.Paragraphs(3).Range
$p3.Find.Execute("Gaude", $true, $false, $false, $false, $false, $true, 1, $false, "Tuyenkar", 2)

# --- 3) Address line: "Address: H No 102 Agapur, Durbhat, Ponda Goa" -> "Address: Marcel, Goa"
$p4 = $d.Paragraphs(4).Range
$p4.Find.Execute("H No 102 Agapur, Durbhat, Ponda Goa", $true, $false, $false, $false, $false, $true, 1, $false, "Marcel, Goa", 2)

# --- 4) Appointed as: "Electrician" -> "Electrical Engineer / Network Planner"
$p5 = $d.Paragraphs(5).Range
$p5.Find.Execute("Electrician", $true, $false, $false, $false, $false, $true, 1, $false, "Electrical Engineer / Network Planner", 2)

# --- 5) Dear line: "Dear Priyesh Gaude," -> "Dear Ganesh Mali,"
$p6 = $d.Paragraphs(6).Range
$p6.Find.Execute(" Priyesh Gaude,", $true, $false, $false, $false, $false, $true, 1, $false, " Ganesh Mali,", 2)

# --- 6) Offer paragraph: "Technical / Service Engineer (Fresher)" -> "Electrical Engineer / Network Planner"
$p7 = $d.Paragraphs(7).Range
$p7.Find.Execute("Technical / Service Engineer (Fresher)", $true, $false, $false, $false, $false, $true, 1, $false, "Electrical Engineer / Network Planner", 2)

# --- 7) Employment begin-on date: "27th January 2024" -> "1st March 2024"
$p10 = $d.Paragraphs(10).Range
$p10.Find.Execute("27", $true, $false, $false, $false, $false, $true, 1, $false, "1", 2)

$p10 = $d.Paragraphs(10).Range
$p10.Find.Execute("th", $true, $false, $false, $false, $false, $true, 1, $false, "st", 2)

$p10 = $d.Paragraphs(10).Range
$p10.Find.Execute("January 2024", $true, $false, $false, $false, $false, $true, 1, $false, "March 2024", 2)

# --- 8) Monthly CTC: "9000" -> "12000"
$p30 = $d.Paragraphs(30).Range
$p30.Find.Execute("9000", $true, $false, $false, $false, $false, $true, 1, $false, "12000", 2)

# --- 9) Remove the _GoBack bookmark currently after the CTC paragraph
$d.Bookmarks("_GoBack").Delete()

# --- 10) Signature line: "(Saish Godkar)" -> "(Ramdas Tuyenkar)"
$p115 = $d.Paragraphs(115).Range
$p115.Find.Execute("Saish", $true, $false, $false, $false, $false, $true, 1, $false, "Ramdas", 2)

$p115b = $d.Paragraphs(115).Range
$p115b.Find.Execute(" Godkar", $true, $false, $false, $false, $false, $true, 1, $false, " Tuyenkar", 2)

# Re-insert the _GoBack bookmark right after "Tuyenkar" (before the closing ")")
$p115c = $d.Paragraphs(115).Range
$p115c.Find.Execute("Tuyenkar", $true, $false, $false, $false, $false, $true)
$bmRange = $p115c.Duplicate
$bmRange.Collapse(0)
$d.Bookmarks.Add("_GoBack", $bmRange)
